# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the target diff: new header cells AD1:AF1 (styled like the other
# header cells) and new data cells AD2:AF42 holding the team's season
# record (78 wins, 84 losses, 0 ties) for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
# Copy the formatting of an existing header cell (bold font, border,
# centered alignment) onto the three new header cells, then set their text.
$ws.Range("A1").Copy($ws.Range("AD1:AF1"))

$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# --- Data ---------------------------------------------------------------
# Every player (rows 2-42) shares the same team season record.
$wins = 78
$losses = 84
$ties = 0

for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value2 = $wins    # column AD
    $ws.Cells.Item($r, 31).Value2 = $losses  # column AE
    $ws.Cells.Item($r, 32).Value2 = $ties    # column AF
}
